$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.845.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.731.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.07"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.04"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.731.99"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.50"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.44"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.354.26"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.720.99"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.789.69"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.05"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.48%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.80"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.25"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.83"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.28%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -11.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.04"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.874.56"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.34"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.77"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.08"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.681.27"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -10.10%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.991"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.58"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.43"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "390.32"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.55"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.41%  "
